$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.07070652293252117, 0.5389307429724479, 0.5493865132072634, 0.7412061205948474, 0.7544083569534742, 23)
    3  = @(0.5139408026494205, 0.8198603642297517, 4.018785774733372, 2.00469094244808, 1.981241167464397, 23)
    4  = @(0.1593057168417592, 1.252584888744187, 7.214834311089271, 2.686044361340533, 2.741577913210072, 23)
    5  = @(0.06397655328083163, 1.169503564658112, 7.127017457775117, 2.669647440726044, 2.72886308780219, 23)
    6  = @(0.1808419384718491, 1.300779946212724, 7.245579745324931, 2.69176145773078, 2.746039682128507, 23)
    7  = @(0.0139830039049362, 1.280086512400622, 7.347913688368553, 2.7107035412174, 2.771588958204612, 23)
    8  = @(0.08558710616064751, 1.381672574370687, 7.474441480762519, 2.733942479417319, 2.794016948153321, 23)
    9  = @(0.007662648018172074, 1.301085017076876, 7.361559590561332, 2.71321941437867, 2.774187187356896, 23)
    10 = @(0.02684109616273803, 1.376516876497369, 7.472336615944717, 2.733557501854445, 2.794858688447044, 23)
    11 = @(-0.129800971266693, 1.190418128082351, 7.078916754308191, 2.660623377013025, 2.717180813652833, 23)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
